$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1): B1 "LF_A" -> "A_V", C1 stays "LF_C", new D1 "LF_FFR"
$ws.Range("B1").Value = "A_V"
$ws.Range("D1").Value = "LF_FFR"

# Give the new header cell the same formatting as the existing header cells
$ws.Range("C1").Copy()
$ws.Range("D1").PasteSpecial(-4122)
$ws.Range("D1").Value = "LF_FFR"

# Row 2 ("params") values
$ws.Range("B2").Value = -0.003229338342638138
$ws.Range("C2").Value = -8.802564485349039
$ws.Range("D2").Value = 0.1641056535089401

# Row 3 ("pvalue") values
$ws.Range("B3").Value = 0.5943035676676558
$ws.Range("C3").Value = 0.01310948737059503
$ws.Range("D3").Value = [double]"2.127116260908224e-10"
